$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(3.227503097112204, 1.2890378437288557, -0.5844423274903815, 2.7566710003242245)"
$ws.Range("C2").Value = "NCT(2.409890526477934, 1.8690899162817805, -1.2119912955382084, 4.479272849903968)"
$ws.Range("D2").Value = "EXN(2.2819188258675425, 0.3401527688869619, 1.6210164165698024)"
$ws.Range("E2").Value = "NIG(1.046255749221365, 0.7312492328597713, 5.611919673700837, 5.213186116445202)"
